$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 2) with new values
$ws.Range("A2").Value = "Kollapudi"
$ws.Range("B2").Value = "Venu"
$ws.Range("C2").Value = "venukollapudi@gmail.com"
$ws.Range("D2").Value = "New York"
$ws.Range("E2").Value = "william Street"
$ws.Range("F2").Value = 10001
$ws.Range("G2").Value = 7013606690

# Set column G width (bestFit-like explicit width) to match new content.
# (10.1666... compensates for the runtime's internal padding so the
#  persisted OOXML <col width="..."> comes out to exactly 11.)
$ws.Columns.Item(7).ColumnWidth = 10.166666666666666

# Move the active selection to G2
$ws.Range("G2").Select()
